$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting used by the existing year labels (column A) down to the
# two new rows, then fill in the values.

$ws.Range("A2").Copy()
$ws.Range("A7").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A7").Value = "2021年"
$ws.Range("B7").Value = 100.4
$ws.Range("C7").Value = 100.8
$ws.Range("D7").Value = 99.5

$ws.Range("A2").Copy()
$ws.Range("A8").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A8").Value = "2022年"
$ws.Range("B8").Value = 100.6
# C8/D8 are intentionally left blank (no data published yet for those columns)

$excel.CutCopyMode = 0
